# Apply "Doing Updates for Financials" edits to the CLHRF yearly financials sheet.
# All changes are numeric-value corrections scattered across the Income
# Statement / Balance Sheet / Cash Flow sections of the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CLHRF")

# Income Statement
$ws.Range("D17").Value = 1000
$ws.Range("D23").Value = 2800
$ws.Range("I23").Value = -1000
$ws.Range("E26").Value = -700
$ws.Range("I26").Value = -1100
$ws.Range("E27").Value = -700
$ws.Range("I27").Value = -1100
$ws.Range("E33").Value = -700
$ws.Range("I33").Value = -1100
$ws.Range("E35").Value = -700
$ws.Range("I35").Value = -1100

# Balance Sheet
$ws.Range("D41").Value = 10700
$ws.Range("F41").Value = 400
$ws.Range("D42").Value = 2700
$ws.Range("D46").Value = 13500
$ws.Range("D47").Value = 2700
$ws.Range("J47").Value = 700
$ws.Range("D48").Value = 100
$ws.Range("E48").Value = 14400
$ws.Range("F48").Value = 14400
$ws.Range("G48").Value = 14300
$ws.Range("H48").Value = 14200
$ws.Range("I48").Value = 14400
$ws.Range("J48").Value = 13400
$ws.Range("D54").Value = 16300
$ws.Range("E54").Value = 14600
$ws.Range("F54").Value = 14900
$ws.Range("G54").Value = 15100
$ws.Range("H54").Value = 14800
$ws.Range("I54").Value = 15100
$ws.Range("J54").Value = 15400
$ws.Range("D60").Value = 2100
$ws.Range("I60").Value = 400
$ws.Range("F62").Value = 900
$ws.Range("I62").Value = 2100
$ws.Range("E66").Value = 1300
$ws.Range("G66").Value = 1500
$ws.Range("J66").Value = 2000

# Cash Flow Statement
$ws.Range("D72").Value = -18900
$ws.Range("E72").Value = -20600
$ws.Range("F72").Value = -19900
$ws.Range("G72").Value = -20100
$ws.Range("H72").Value = -20300
$ws.Range("I72").Value = -22400
$ws.Range("J72").Value = -21600
$ws.Range("D76").Value = 14200
$ws.Range("E76").Value = 13300
$ws.Range("F76").Value = 13900
$ws.Range("G76").Value = 13600
$ws.Range("H76").Value = 12700
$ws.Range("I76").Value = 12600
$ws.Range("J76").Value = 13400

# Balance Sheet (second block, below Cash Flow)
$ws.Range("E81").Value = -700
$ws.Range("I81").Value = -1100
$ws.Range("E89").Value = -300
$ws.Range("J89").Value = -500
$ws.Range("D94").Value = 12200
$ws.Range("D102").Value = 10600
$ws.Range("J102").Value = -1000
